$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1970
$ws.Range("I43").Value = 1970
$ws.Range("K43").Value = 1970
$ws.Range("M43").Value = -1901
$ws.Range("H80").Value = 819
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 819
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 2457
$ws.Range("N80").Value = -4453
$ws.Range("H83").Value = 819
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 819
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 7371
$ws.Range("N83").Value = -17355
$ws.Range("H101").Value = 20003764
$ws.Range("I101").Value = 33339114
$ws.Range("J101").Value = 737.5
$ws.Range("K101").Value = 100017342
$ws.Range("L101").Value = 2212.5
$ws.Range("M101").Value = -100015720
$ws.Range("N101").Value = -5456.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0
$ws.Range("H132").Value = 2462.9092
$ws.Range("I132").Value = 2632.6667
$ws.Range("K132").Value = 7898.000100000001
$ws.Range("M132").Value = -5368.000100000001
$ws.Range("H137").Value = 4524.5
$ws.Range("I137").Value = 2412.125
$ws.Range("J137").Value = 8749.25
$ws.Range("K137").Value = 7236.375
$ws.Range("L137").Value = 26247.75
$ws.Range("M137").Value = -4686.375
$ws.Range("N137").Value = -31347.75
$ws.Range("H141").Value = 3999.5
$ws.Range("I141").Value = 3999.5
$ws.Range("K141").Value = 11998.5
$ws.Range("M141").Value = -6818.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3966.6667
$ws.Range("H111").Value = 8000
$ws.Range("J111").Value = 8000
$ws.Range("L111").Value = 8000
$ws.Range("N111").Value = -16180
$ws.Range("H116").Value = 3966.6667
$ws.Range("H132").Value = 899.4286
$ws.Range("I132").Value = 774.75
$ws.Range("J132").Value = 1065.6666
$ws.Range("K132").Value = 2324.25
$ws.Range("L132").Value = 3196.9998
$ws.Range("M132").Value = 205.75
$ws.Range("N132").Value = -8256.9998
$ws.Range("H135").Value = 1000000
$ws.Range("J135").Value = 1000000
$ws.Range("L135").Value = 1000000
$ws.Range("N135").Value = -1010140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3966.6667
$ws.Range("H86").Value = 1891.3572
$ws.Range("I86").Value = 1768
$ws.Range("J86").Value = 2199.75
$ws.Range("K86").Value = 1768
$ws.Range("L86").Value = 2199.75
$ws.Range("M86").Value = -645
$ws.Range("N86").Value = -4445.75
$ws.Range("H89").Value = 1891.3572
$ws.Range("I89").Value = 1768
$ws.Range("J89").Value = 2199.75
$ws.Range("K89").Value = 8840
$ws.Range("L89").Value = 10998.75
$ws.Range("M89").Value = -3224
$ws.Range("N89").Value = -22230.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3611.3333
$ws.Range("I16").Value = 4688.25
$ws.Range("J16").Value = 2749.8
$ws.Range("K16").Value = 4688.25
$ws.Range("L16").Value = 2749.8
$ws.Range("M16").Value = -4401.25
$ws.Range("N16").Value = -3323.8
$ws.Range("H31").Value = 1969.875
$ws.Range("I31").Value = 1969.875
$ws.Range("K31").Value = 1969.875
$ws.Range("M31").Value = -1674.875
$ws.Range("H34").Value = 1969.875
$ws.Range("I34").Value = 1969.875
$ws.Range("K34").Value = 1969.875
$ws.Range("M34").Value = -1767.875
$ws.Range("H99").Value = 1896.3572
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 2137.25
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2137.25
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -5133.25
$ws.Range("H113").Value = 3611.3333
$ws.Range("I113").Value = 4688.25
$ws.Range("J113").Value = 2749.8
$ws.Range("K113").Value = 4688.25
$ws.Range("L113").Value = 2749.8
$ws.Range("M113").Value = -2518.25
$ws.Range("N113").Value = -7089.8
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 925
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 2775
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -325
$ws.Range("N122").Value = -8800
$ws.Range("H126").Value = 1896.3572
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2137.25
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 6411.75
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -11351.75
$ws.Range("H132").Value = 3993.75
$ws.Range("I132").Value = 3992.6667
$ws.Range("K132").Value = 11978.0001
$ws.Range("M132").Value = -9448.000100000001
$ws.Range("H134").Value = 2574.1667
$ws.Range("J134").Value = 1999.5
$ws.Range("L134").Value = 5998.5
$ws.Range("N134").Value = -11068.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2697.3157
$ws.Range("J22").Value = 2736.0557
$ws.Range("L22").Value = 8208.167099999999
$ws.Range("N22").Value = -8546.167099999999
$ws.Range("H23").Value = 42857296
$ws.Range("I23").Value = 75000090
$ws.Range("J23").Value = 241.66667
$ws.Range("K23").Value = 225000270
$ws.Range("L23").Value = 725.00001
$ws.Range("M23").Value = -225000035
$ws.Range("N23").Value = -1195.00001
$ws.Range("H27").Value = 2697.3157
$ws.Range("J27").Value = 2736.0557
$ws.Range("L27").Value = 8208.167099999999
$ws.Range("N27").Value = -8412.167099999999
$ws.Range("H63").Value = 16786.666
$ws.Range("J63").Value = 1900
$ws.Range("L63").Value = 5700
$ws.Range("N63").Value = -7198
$ws.Range("H66").Value = 16786.666
$ws.Range("J66").Value = 1900
$ws.Range("L66").Value = 17100
$ws.Range("N66").Value = -24588
$ws.Range("H81").Value = 2400
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 2400
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H106").Value = 16482.834
$ws.Range("I106").Value = 9449.5
$ws.Range("K106").Value = 28348.5
$ws.Range("M106").Value = -27402.5
$ws.Range("H118").Value = 1500
$ws.Range("I118").Value = 1500
$ws.Range("K118").Value = 4500
$ws.Range("M118").Value = -3257
$ws.Range("H131").Value = 590075.44
$ws.Range("I131").Value = 1175.375
$ws.Range("J131").Value = 1113542.1
$ws.Range("K131").Value = 3526.125
$ws.Range("L131").Value = 3340626.3
$ws.Range("M131").Value = 1513.875
$ws.Range("N131").Value = -3350706.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 15040000
$ws.Range("J33").Value = 15040000
$ws.Range("L33").Value = 15040000
$ws.Range("N33").Value = -15040504
$ws.Range("H126").Value = 5933.3335
$ws.Range("I126").Value = 800
$ws.Range("K126").Value = 2400
$ws.Range("M126").Value = 70

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1591.7858
$ws.Range("I46").Value = 1996
$ws.Range("J46").Value = 1481.5454
$ws.Range("K46").Value = 1996
$ws.Range("L46").Value = 1481.5454
$ws.Range("M46").Value = -1808
$ws.Range("N46").Value = -1857.5454
$ws.Range("H55").Value = 694.6087
$ws.Range("I55").Value = 692.5
$ws.Range("J55").Value = 696.9091
$ws.Range("K55").Value = 692.5
$ws.Range("L55").Value = 696.9091
$ws.Range("M55").Value = -519.5
$ws.Range("N55").Value = -1042.9091
$ws.Range("H122").Value = 6610.5386
$ws.Range("I122").Value = 5437
$ws.Range("K122").Value = 16311
$ws.Range("M122").Value = -13861
$ws.Range("H132").Value = 3564.5
$ws.Range("I132").Value = 3631
$ws.Range("K132").Value = 10893
$ws.Range("M132").Value = -8363
$ws.Range("H136").Value = 2272.5715
$ws.Range("J136").Value = 3002
$ws.Range("L136").Value = 9006
$ws.Range("N136").Value = -14106

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 50005
$ws.Range("J11").Value = 50005
$ws.Range("L11").Value = 50005
$ws.Range("N11").Value = -50289
$ws.Range("H126").Value = 3480.6
$ws.Range("I126").Value = 2004
$ws.Range("K126").Value = 6012
$ws.Range("M126").Value = -3542
$ws.Range("H132").Value = 1984
$ws.Range("J132").Value = 699
$ws.Range("L132").Value = 2097
$ws.Range("N132").Value = -7157
